$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table used to start one row down / one column right (row 2, column B)
# with an empty spacer column A and an empty spacer row 1. The edit removes
# those spacers: the header + 5 lockdown rows move up by one row and left
# by one column (B:D -> A:C, rows 2:7 -> rows 1:6).

$ws.Range("B2:D2").Copy($ws.Range("A1:C1"))
$ws.Range("B3:D3").Copy($ws.Range("A2:C2"))
$ws.Range("B4:D4").Copy($ws.Range("A3:C3"))
$ws.Range("B5:D5").Copy($ws.Range("A4:C4"))
$ws.Range("B6:D6").Copy($ws.Range("A5:C5"))
$ws.Range("B7:D7").Copy($ws.Range("A6:C6"))

# Wipe what's left of the old table in column D and the old row 7's B:C
# (row 7 col A already carries the correct blank-spacer style, untouched).
$ws.Range("D2:D7").Clear()
$ws.Range("B7:C7").Clear()

# Match the saved cursor position from the edit.
$ws.Range("B13").Select() | Out-Null
